$d = $word.ActiveDocument

# --- 1) OC-1 section: "Krydsreferencer" body paragraph (5th paragraph)
#        "opretOrganisation" -> "Opret organisation"
$r = $d.Paragraphs.Item(5).Range
$r.Find.Execute("opretOrganisation", $false, $false, $false, $false, $false, $true, 1, $false, "Opret organisation", 2)

# --- 2) OC-2 section: "Krydsreferencer" body paragraph (20th paragraph)
#        "opretOrganisation" -> "Opret organisation"
$r = $d.Paragraphs.Item(20).Range
$r.Find.Execute("opretOrganisation", $false, $false, $false, $false, $false, $true, 1, $false, "Opret organisation", 2)

# --- 3) OC-2 section: "Slutbetingelser" body paragraph (26th paragraph)
#        "Præsenter de angivne oplysninger. " -> two new paragraphs describing the
#        created Organisation instance and its presentation.
$r = $d.Paragraphs.Item(26).Range
$r.MoveEnd(1, -1)
$newText = "En instans o af Organisation er blevet skabt." + [char]11 + `
    "o.CVR er sat til CVR" + [char]11 + `
    "o.navn er sat til navn" + [char]11 + `
    "o.email er sat til " + "email" + [char]11 + `
    "o.adresse er sat til adresse" + [char]11 + `
    "o.tlf er sat til " + "tlf" + [char]13 + `
    "o er blevet præsenteret"
$r.Text = $newText

# --- 4) OC-3 section: "Krydsreferencer" body paragraph.
#        Was paragraph 34 before the insertion above added a paragraph; now 35.
#        "opretOrganisation." -> "Opret organisation."
$r = $d.Paragraphs.Item(35).Range
$r.Find.Execute("opretOrganisation", $false, $false, $false, $false, $false, $true, 1, $false, "Opret organisation", 2)
